# Generate Report for Handoff
# The file "abf1d66b-f12b-48e9-b67f-01f10da0c3c9.md" has finished translation
# and is now ready for handoff: update its status on every sheet, and stamp
# fresh "Latest Handoff Datetime" values on the locale sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: update the per-locale status columns (B3, C3) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: update status + handoff datetime for abf1d66b...md (row 3) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-03-09 12:12:59"

# --- de-de sheet: update status + handoff datetime for abf1d66b...md (row 3) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-03-09 12:13:05"
